$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "semana 31 de 2025" update
#   1) Insert a brand-new institution row at row 39 (SALUD PYP SAS), which
#      pushes every following row down by one (old row 39 -> new row 40, ...,
#      old row 57 -> new row 58).
#   2) Add a new week-31 column (AH) with a header cell and one data value
#      per institution row.
# ---------------------------------------------------------------------------

# 1) Insert the new row for "SALUD PYP SAS" at position 39.
$ws.Rows.Item(39).Insert()

$ws.Range("A39:B39").NumberFormat = "@"
$ws.Range("A39").Value = "6600102288"
$ws.Range("B39").Value = "01"
$ws.Range("A39:B39").Style = "Normal"

$ws.Range("C39").Value = "SALUD PYP SAS"

# 2) New week-31 header in column AH, formatted like the other week headers.
$ws.Range("AH1").NumberFormat = "@"
$ws.Range("AH1").Value = "31"
$ws.Range("AG1").Copy()
$ws.Range("AH1").PasteSpecial(-4122)

# 3) Week-31 counts for every institution row that reported a value.
$ws.Range("AH2").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AH6").Value = 2
$ws.Range("AH7").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AH10").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AH12").Value = 0
$ws.Range("AH13").Value = 0
$ws.Range("AH14").Value = 0
$ws.Range("AH15").Value = 0
$ws.Range("AH16").Value = 0
$ws.Range("AH17").Value = 0
$ws.Range("AH23").Value = 0
$ws.Range("AH24").Value = 0
$ws.Range("AH25").Value = 0
$ws.Range("AH28").Value = 103
$ws.Range("AH29").Value = 0
$ws.Range("AH30").Value = 4
$ws.Range("AH31").Value = 0
$ws.Range("AH32").Value = 0
$ws.Range("AH34").Value = 0
$ws.Range("AH35").Value = 2
$ws.Range("AH36").Value = 0
$ws.Range("AH37").Value = 0
$ws.Range("AH38").Value = 0
$ws.Range("AH39").Value = 0
$ws.Range("AH40").Value = 0
$ws.Range("AH41").Value = 0
$ws.Range("AH42").Value = 0
$ws.Range("AH43").Value = 0
$ws.Range("AH45").Value = 0
$ws.Range("AH46").Value = 0
$ws.Range("AH47").Value = 0
$ws.Range("AH48").Value = 0
$ws.Range("AH49").Value = 0
$ws.Range("AH50").Value = 0
$ws.Range("AH51").Value = 0
$ws.Range("AH53").Value = 0
$ws.Range("AH54").Value = 0
$ws.Range("AH55").Value = 0
$ws.Range("AH56").Value = 0
$ws.Range("AH57").Value = 0
$ws.Range("AH58").Value = 0
